$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'303.64"
$ws.Cells.Item(2, 5).Value = "'5.53%"
$ws.Cells.Item(3, 4).Value = "'35.02"
$ws.Cells.Item(3, 5).Value = "'12.62%"
$ws.Cells.Item(4, 4).Value = "'5.155"
$ws.Cells.Item(4, 5).Value = "'4.83%"
$ws.Cells.Item(5, 4).Value = "'0.07755"
$ws.Cells.Item(5, 5).Value = "'5.94%"
$ws.Cells.Item(6, 4).Value = "'2.296"
$ws.Cells.Item(6, 5).Value = "'3.06%"
$ws.Cells.Item(7, 4).Value = "'8.011"
$ws.Cells.Item(7, 5).Value = "'3.59%"
$ws.Cells.Item(8, 4).Value = "'3.985"
$ws.Cells.Item(8, 5).Value = "'7.12%"
$ws.Cells.Item(9, 4).Value = "'0.9249"
$ws.Cells.Item(9, 5).Value = "'2.27%"
$ws.Cells.Item(10, 4).Value = "'0.1013"
$ws.Cells.Item(10, 5).Value = "'10.80%"
$ws.Cells.Item(11, 4).Value = "'0.1838"
$ws.Cells.Item(11, 5).Value = "'8.94%"
$ws.Cells.Item(12, 4).Value = "'0.08577"
$ws.Cells.Item(12, 5).Value = "'4.71%"
$ws.Cells.Item(13, 4).Value = "'0.03401"
$ws.Cells.Item(13, 5).Value = "'8.90%"
$ws.Cells.Item(14, 4).Value = "'0.09885"
$ws.Cells.Item(14, 5).Value = "'-0.41%"
$ws.Cells.Item(15, 4).Value = "'0.001478"
$ws.Cells.Item(15, 5).Value = "'-1.32%"
$ws.Cells.Item(16, 4).Value = "'0.005766"
$ws.Cells.Item(16, 5).Value = "'0.88%"
$ws.Cells.Item(17, 4).Value = "'3.504"
$ws.Cells.Item(17, 5).Value = "'-0.90%"
$ws.Cells.Item(18, 4).Value = "'2.108"
$ws.Cells.Item(18, 5).Value = "'1.29%"
$ws.Cells.Item(19, 4).Value = "'0.3441"
$ws.Cells.Item(19, 5).Value = "'3.33%"
$ws.Cells.Item(20, 4).Value = "'0.1305"
$ws.Cells.Item(20, 5).Value = "'0.51%"
$ws.Cells.Item(21, 4).Value = "'4.379"
$ws.Cells.Item(21, 5).Value = "'4.24%"
$ws.Cells.Item(22, 4).Value = "'0.2299"
$ws.Cells.Item(22, 5).Value = "'9.35%"
$ws.Cells.Item(23, 4).Value = "'0.04627"
$ws.Cells.Item(23, 5).Value = "'2.78%"
$ws.Cells.Item(24, 4).Value = "'0.001221"
$ws.Cells.Item(24, 5).Value = "'0.76%"
$ws.Cells.Item(25, 4).Value = "'0.004444"
$ws.Cells.Item(25, 5).Value = "'6.91%"
$ws.Cells.Item(26, 4).Value = "'0.0001309"
$ws.Cells.Item(26, 5).Value = "'0.56%"
$ws.Cells.Item(27, 4).Value = "'0.0003396"
$ws.Cells.Item(27, 5).Value = "'-0.06%"
$ws.Cells.Item(39, 4).Value = "'0.01796"
$ws.Cells.Item(39, 5).Value = "'14.33%"
$ws.Cells.Item(40, 4).Value = "'0.04735"
$ws.Cells.Item(40, 5).Value = "'6.71%"
$ws.Cells.Item(41, 4).Value = "'0.007521"
$ws.Cells.Item(41, 5).Value = "'2.68%"
$ws.Cells.Item(42, 4).Value = "'0.1403"
$ws.Cells.Item(42, 5).Value = "'5.70%"
$ws.Cells.Item(43, 4).Value = "'0.007050"
$ws.Cells.Item(43, 5).Value = "'-26.34%"
$ws.Cells.Item(44, 4).Value = "'0.002233"
$ws.Cells.Item(44, 5).Value = "'-2.61%"
$ws.Cells.Item(45, 4).Value = "'0.009240"
$ws.Cells.Item(45, 5).Value = "'10.81%"
$ws.Cells.Item(46, 4).Value = "'0.00005767"
$ws.Cells.Item(46, 5).Value = "'-5.61%"
$ws.Cells.Item(47, 4).Value = "'0.00000000755"
$ws.Cells.Item(47, 5).Value = "'0.51%"
$ws.Cells.Item(48, 5).Value = "'19.99%"
$ws.Cells.Item(49, 4).Value = "'0.002686"
$ws.Cells.Item(49, 5).Value = "'34.18%"
$ws.Cells.Item(50, 4).Value = "'0.00002113"
$ws.Cells.Item(50, 5).Value = "'0.51%"
$ws.Cells.Item(51, 4).Value = "'0.0002012"
$ws.Cells.Item(51, 5).Value = "'0.51%"
